$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 25; $row -le 81; $row++) {
    $ws.Cells.Item($row, 19).Value = "https://orcid.org/0000-0003-2195-3997"
}

$ws.Cells.Item(25, 21).Value = "2023-08-26"
